$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
$ws.Range("D2").Value = "67.832.58"
$ws.Range("D3").Value = "3.543.28"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "616.47"
$ws.Range("D7").Value = "3.539.57"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000224"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.23"
$ws.Range("D15").Value = "4.141.47"
$ws.Range("D16").Value = "3.541.97"
$ws.Range("D17").Value = "67.797.62"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.38"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.52"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "453.87"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.10"
$ws.Range("D25").Value = "3.683.83"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.56"
$ws.Range("D37").Value = "3.540.71"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.14"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0881"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.50"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.68"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -3.85%  "
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("E7").Value = "  -3.40%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("E15").Value = "  -3.43%  "
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -3.66%  "
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("E46").Value = "  +10.11%  "
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("E51").Value = "  -2.93%  "
